$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style from H1 (bold, bordered, centered) onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows for columns I and J
$data = @{
    2 = @(1, 5)
    3 = @(1, 6)
    4 = @(1, 4)
    5 = @(1, 5)
    6 = @(1, 5)
    7 = @(1, 6)
    8 = @(8, 8)
    9 = @(9, 9)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
